$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns O1 and P1
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Data rows 2-8, columns O and P
$ws.Range("O2").Value = -0.8227941827338482
$ws.Range("P2").Value = -0.6044233549880282

$ws.Range("O3").Value = -0.4550981943703185
$ws.Range("P3").Value = -0.3910084373000415

$ws.Range("O4").Value = 0.03520707643816495
$ws.Range("P4").Value = -0.01096204100079212

$ws.Range("O5").Value = 0.4132367660373626
$ws.Range("P5").Value = 0.373442009009972

$ws.Range("O6").Value = -0.2930756934552439
$ws.Range("P6").Value = -0.2838677905904577

$ws.Range("O7").Value = -0.1523061062630485
$ws.Range("P7").Value = -0.1516627028234634

$ws.Range("O8").Value = -0.4283265539945751
$ws.Range("P8").Value = -0.4197994409469847
